$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Rename the AddressBook-related identifiers to TravelBuddy, per the commit:
#  "Modify "addressbook", "address book", "person" and "persons" to
#   "travelbuddy", "place" or "places"

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $tr = $sh.TextFrame.TextRange
        $null = $tr.Replace("handleAddresssBookChangedEvent", "handleTravelBuddyChangedEvent", 0, $false, $false)
        $null = $tr.Replace("AddressBookChangedEvent", "TravelBuddyChangedEvent", 0, $false, $false)
    }
}
